# Insert a new weekly price record for "Ciboulette" (Feria Lagunitas de Puerto
# Montt) immediately before the current row 69. Excel shifts every existing
# row at/after 69 down by one (old row 69 -> 70, ..., old row 109 -> 110) and
# extends the used range to R110; we then populate the freshly inserted row
# 69 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 69:109 down one row, opening up a blank row 69.
$ws.Rows(69).Insert()

# Fill the new row 69 with the new data point.
$ws.Cells.Item(69, 1).Value  = 4
$ws.Cells.Item(69, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(69, 3).Value  = "Los Lagos"
$ws.Cells.Item(69, 4).Value  = 44452
$ws.Cells.Item(69, 5).Value  = 10
$ws.Cells.Item(69, 6).Value  = 100112039
$ws.Cells.Item(69, 7).Value  = "Ciboulette"
$ws.Cells.Item(69, 8).Value  = "Sin especificar"
$ws.Cells.Item(69, 9).Value  = "Primera"
$ws.Cells.Item(69, 10).Value = 80
$ws.Cells.Item(69, 11).Value = 4500
$ws.Cells.Item(69, 12).Value = 4500
$ws.Cells.Item(69, 13).Value = 4500
$ws.Cells.Item(69, 14).Value = "`$/docena de atados"
$ws.Cells.Item(69, 15).Value = "Región Metropolitana"
$ws.Cells.Item(69, 16).Value = 1500
$ws.Cells.Item(69, 17).Value = 3
$ws.Cells.Item(69, 18).Value = "Hortaliza"
